$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 202-203. This pushes the existing rows 202-246
# down to 204-248, matching the rest of the diff (which is a pure shift
# of that block by two rows with no content changes beyond the shift).
$ws.Rows("202:203").Insert()

# New row 202 data
$ws.Cells.Item(202, 1).Value = 5
$ws.Cells.Item(202, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(202, 3).Value = "Maule"
$ws.Cells.Item(202, 4).Value = 44559
$ws.Cells.Item(202, 5).Value = 7
$ws.Cells.Item(202, 6).Value = "Fruta"
$ws.Cells.Item(202, 7).Value = 100103
$ws.Cells.Item(202, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(202, 9).Value = 100103004
$ws.Cells.Item(202, 10).Value = "Durazno"
$ws.Cells.Item(202, 11).Value = "Kurakata"
$ws.Cells.Item(202, 12).Value = "Primera"
$ws.Cells.Item(202, 13).Value = 230
$ws.Cells.Item(202, 14).Value = 10000
$ws.Cells.Item(202, 15).Value = 10000
$ws.Cells.Item(202, 16).Value = 10000
$ws.Cells.Item(202, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(202, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(202, 19).Value = 667
$ws.Cells.Item(202, 20).Value = 15

# New row 203 data
$ws.Cells.Item(203, 1).Value = 5
$ws.Cells.Item(203, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(203, 3).Value = "Maule"
$ws.Cells.Item(203, 4).Value = 44559
$ws.Cells.Item(203, 5).Value = 7
$ws.Cells.Item(203, 6).Value = "Fruta"
$ws.Cells.Item(203, 7).Value = 100103
$ws.Cells.Item(203, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(203, 9).Value = 100103004
$ws.Cells.Item(203, 10).Value = "Durazno"
$ws.Cells.Item(203, 11).Value = "Toscana"
$ws.Cells.Item(203, 12).Value = "Especial"
$ws.Cells.Item(203, 13).Value = 180
$ws.Cells.Item(203, 14).Value = 12000
$ws.Cells.Item(203, 15).Value = 12000
$ws.Cells.Item(203, 16).Value = 12000
$ws.Cells.Item(203, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(203, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(203, 19).Value = 800
$ws.Cells.Item(203, 20).Value = 15
